$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$refStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "68.950.73"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "3.934.35"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'488.21"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "'146.72"
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("D11").Value = "'0.0000344"
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = "  -4.77%  "
$ws.Range("D12").Value = "'43.06"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "'10.47"
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "4.555.52"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "3.944.91"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "'14.24"
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D18").Value = "'19.96"
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "'1.17"
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("D20").Value = "68.927.92"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").Value = "'437.13"
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("E22").Value = "  +3.11%  "
$ws.Range("D23").Value = "'14.59"
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "'12.48"
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = "  +17.15%  "
$ws.Range("D25").Value = "'89.43"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "'3.73"
$ws.Range("D26").Style = $refStyle
$ws.Range("E26").Value = "  +3.48%  "
$ws.Range("D27").Value = "'11.09"
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").Value = "'37.18"
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("E29").Value = "  -3.62%  "
$ws.Range("D30").Value = "'710.42"
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").Value = "'13.54"
$ws.Range("D31").Style = $refStyle
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "'0.132"
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "'0.480"
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = "  +31.07%  "
$ws.Range("D35").Value = "0.0₃0890"
$ws.Range("E35").Value = "  -5.38%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'61.68"
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'6.06"
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = "  +6.85%  "
$ws.Range("D38").Value = "'40.82"
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'2.97"
$ws.Range("D42").Style = $refStyle
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("D43").Value = "'0.0490"
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0368"
$ws.Range("E47").Value = "  +10.67%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.35"
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = "  +6.44%  "
$ws.Range("D49").Value = "'3.03"
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = "  +6.31%  "
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E51").Value = "  -2.74%  "
